# North_Valley_Data.xlsx — "data reading is complete"
#
# 1) Column A (rows 3-20) held the sample year as a lookup into the shared
#    string table (e.g. "2010"). Data reading is complete, so the years are
#    now entered as real numbers instead of text.
# 2) A new aggregate/mix sample (row 21) is appended below the existing
#    20 rows, with a trailing blank, bold-formatted spacer row (row 22).
# 3) The active selection moves down to D22, where the next entry would go.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Column A: years as numbers instead of text -------------------------
$years = @{
    3  = 2010
    4  = 2012
    5  = 2016
    6  = 2018
    7  = 2002
    8  = 2003
    9  = 2004
    10 = 2005
    11 = 2006
    12 = 2007
    13 = 2008
    14 = 2009
    15 = 2011
    16 = 2013
    17 = 2014
    18 = 2017
    19 = 2015
    20 = 2001
}
foreach ($row in $years.Keys) {
    $ws.Range("A$row").Value = $years[$row]
}

# --- 2) New row 21: aggregated "Mix" sample ---------------------------------
$boldCells = @{
    "A21"  = 3000
    "E21"  = $null
    "F21"  = 5.03
    "G21"  = 427.3
    "H21"  = 0.9
    "I21"  = $null
    "J21"  = $null
    "K21"  = 3.68
    "L21"  = 2.9
    "M21"  = 0.81
    "N21"  = $null
    "O21"  = 2.1
    "P21"  = 12.35
    "Q21"  = $null
    "R21"  = $null
    "S21"  = $null
    "T21"  = 1.2
    "U21"  = $null
    "V21"  = 0.3
    "W21"  = 0.6
    "AK21" = 9
    "AL21" = $null
    "AP21" = 59
    "AQ21" = 13
    "AR21" = 17
}
foreach ($addr in $boldCells.Keys) {
    $cell = $ws.Range($addr)
    if ($null -ne $boldCells[$addr]) {
        $cell.Value = $boldCells[$addr]
    }
    $cell.Font.Bold = $true
}

# hlorit (AM21) carries the mix's dimensionless "2" code, stored as text
# (routed through a text formula + value-bake so it lands as a plain shared
# string instead of being auto-coerced to the number 2)
$ws.Range("AM21").Formula = '="2"'
$ws.Range("AM21").Value = $ws.Range("AM21").Value
$ws.Range("AM21").Font.Bold = $true

# Source (D21) is plain, unbolded text
$ws.Range("D21").Value = "Mix"

# --- 3) New row 22: blank bold spacer row under the same columns -----------
$spacerCols = @("A","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","AK","AL","AM","AP","AQ","AR")
foreach ($col in $spacerCols) {
    $ws.Range("$col" + "22").Font.Bold = $true
}

# --- 4) Move the active selection to D22 ------------------------------------
$ws.Range("D22").Select()
